$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (dLacI): bmin 1E-3 -> 1E-4, bmax 10 -> 1 ---
$ws.Range("B4").Value = 0.0001
$ws.Range("C4").Value = 1

# --- Row 5 (dCit): bmin 7E-3 -> 1E-4, bmax 8E-3 -> 1, estimate "no" -> "yes" ---
$ws.Range("B5").Value = 0.0001
$ws.Range("C5").Value = 1
$ws.Range("F5").Value = "yes"

# --- Row 9 (LacI_rep_Cit): bmin 1E-3 -> 1E-5 ---
$ws.Range("B9").Value = 0.00001

# --- Row 10 (LacI_rep_Cit_W220F): bmin 1E-3 -> 1E-5 ---
$ws.Range("B10").Value = 0.00001

# --- Row 16 (new): P_4Lacn_LacI ---
$ws.Range("A16").Value = "P_4Lacn_LacI"
$ws.Range("B16").Value = 0.01
$ws.Range("C16").Value = 100
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 98
$ws.Range("F16").Value = "yes"
$ws.Range("G16").Value = "k_{LacI_W220F_Q60G_T167A}"
$ws.Range("G16").Font.Size = 11

# --- Row 17 (new): P_4Lacn_LacI_L ---
$ws.Range("A17").Value = "P_4Lacn_LacI_L"
$ws.Range("B17").Value = 0.00001
$ws.Range("C17").Value = 0.01
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0.0003
$ws.Range("F17").Value = "yes"
$ws.Range("G17").Value = "kL_W220F_Q60G_T167A"
$ws.Range("G17").Font.Size = 11

# --- Row 18 (new): LacI_rep ---
$ws.Range("A18").Value = "LacI_rep"
$ws.Range("B18").Value = 0.00001
$ws.Range("C18").Value = 100
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = "yes"
$ws.Range("G18").Value = "theta_{LacI_W220F_Q60G_T167A}"

# --- Selection moved from C7 to H21 ---
$null = $ws.Range("H21").Select()
